$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out the old second data row's extra columns (B,C,D,F,G,H,I) and
# replace the SONG_ID values with the new placeholder unique ids, adding a
# third row so the available ids start from row 2 ("1" relative index).
$ws.Range("B2:I2").Clear()

$ws.Range("A2").Value = "pkmdkm"
$ws.Range("A3").Value = "gdmsen"

$ws.Range("A4").Select()
